$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells: Wins / Losses / Ties in columns AD, AE, AF (row 1)
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Match the header formatting used by the other header cells (e.g. AC1) -
# bold font + thin border + centered/top alignment. PasteSpecial(xlPasteFormats)
# only copies formatting, leaving the values set above untouched.
$headerSrc = $ws.Range("AC1")
$headerDst = $ws.Range("AD1:AF1")
$headerSrc.Copy()
$headerDst.PasteSpecial(-4122)

# Fill team record (Wins=75, Losses=87, Ties=0) for every data row (2-36)
$lastRow = 36
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 30).Value = 75
    $ws.Cells.Item($r, 31).Value = 87
    $ws.Cells.Item($r, 32).Value = 0
}
